$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_SprintBacklog")

# --- Column B width update (~23.875 chars; engine quantizes to the Normal-font pixel grid) ---
$ws.Columns.Item(2).ColumnWidth = 23.16

# --- Task rows: update Id numbers, Status, and "pozostało" (remaining hours) column E ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "IN PROGRESS"
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "IN PROGRESS"
$ws.Range("E6").Value = 1

$ws.Range("A7").Value = 3
$ws.Range("E7").Value = 2.5

$ws.Range("A8").Value = 4
$ws.Range("E8").Value = 2.5

$ws.Range("A9").Value = 5
$ws.Range("E9").Value = 2

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "IN PROGRESS"
$ws.Range("E10").Value = 0.5

$ws.Range("A11").Value = 7
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 20

$ws.Range("A12").Value = 8
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 20

$ws.Range("A13").Value = 9
$ws.Range("E13").Value = 2

$ws.Range("A14").Value = 10
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 10

$ws.Range("A15").Value = 11
$ws.Range("E15").Value = 1

$ws.Range("A16").Value = 12
$ws.Range("E16").Value = 0.5

$ws.Range("A17").Value = 13
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 9

$ws.Range("A18").Value = 14
$ws.Range("E18").Value = 2

$ws.Range("A19").Value = 15
$ws.Range("E19").Value = 4

$ws.Range("A20").Value = 16
$ws.Range("E20").Value = 1

# --- Totals row (values set further below, after new strings are registered) ---
$ws.Range("D21").Formula = "=SUM(D5:D20)"
$ws.Range("E21").Formula = "=SUM(E5:E20)"

# --- Burndown data table (rows 29-65) ---
$ws.Range("A29").Value = "data"
$ws.Range("B29").Value = "pozostało [h]"

# --- Totals row label (added last so "SUMA" lands after "data"/"pozostało [h]" in the shared strings table) ---
$ws.Range("A21").Value = "SUMA"

$ws.Range("A30").Value = 41273
$ws.Range("B30").Value = 80.5
$ws.Range("A31").Value = 41274
$ws.Range("B31").Value = 79

$r = 32
for ($serial = 41275; $serial -le 41308; $serial++) {
    $ws.Cells.Item($r, 1).Value = $serial
    $r = $r + 1
}

$ws.Range("A30:A65").NumberFormat = "yyyy/mm/dd"

# --- Burndown chart ---
$co = $ws.ChartObjects().Add(1143000, 4400000, 2360000, 2600000)
$chart = $co.Chart
$chart.ChartType = 4
$chart.SeriesCollection.NewSeries()
$ser = $chart.SeriesCollection(1)
$ser.Name = "Sprint1"
$ser.XValues = $ws.Range("A30:A65")
$ser.Values = $ws.Range("B30:B65")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "burndown"

$chart.HasLegend = $true
$chart.Legend.Position = -4152

$valAx = $chart.Axes(2)
$valAx.MinimumScale = 0
$valAx.MaximumScale = 81

# --- Selection / view ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B32").Select()

$wb.Save()
